$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (K2:T2)
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4511956666666667
$ws.Range("N2").Value = 1.353587
$ws.Range("O2").Value = 0.7507457057443468
$ws.Range("P2").Value = 0.750745705744347
$ws.Range("Q2").Value = 0.1465889601433334
$ws.Range("R2").Value = 1.31930064129
$ws.Range("S2").Value = 0.7507457057443468
$ws.Range("T2").Value = 0.750745705744347

# Row 3 updates (K3:T3)
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.149801
$ws.Range("N3").Value = 0.449403
$ws.Range("O3").Value = 0.2492542942556531
$ws.Range("P3").Value = 0.2492542942556531
$ws.Range("Q3").Value = 0.04866884689
$ws.Range("R3").Value = 0.43801962201
$ws.Range("S3").Value = 0.2492542942556531
$ws.Range("T3").Value = 0.2492542942556531
